# Updates the "delete 1" mini sequence-diagram (shapes 12 & 14 on slide 1)
# to the new "delete-patient n/John Doe" command, matching the commit that
# added several more sequence diagrams to the developer guide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "TextBox 23" (the actor-side "delete 1" label) -----------------
$sh1 = $s.Shapes.Item("TextBox 23")

# Reposition / resize (EMU targets: off 434250,1328746 ext 1072028,430887).
# Point values below are chosen so the COM layer's single-precision/EMU
# rounding lands exactly on the target EMU values.
$sh1.Left = 34.19291498582677
$sh1.Top = 104.62566769133858
$sh1.Width = 84.41165554330708
$sh1.Height = 33.92811023622047

$sh1.TextFrame.TextRange.Text = "delete-patient n/John Doe"

# --- Shape "TextBox 25" (the "execute(...)" label) -------------------------
$sh2 = $s.Shapes.Item("TextBox 25")

# Reposition / resize (EMU targets: off 1981678,1464043 ext 1760984,430887).
$sh2.Left = 156.03763599527562
$sh2.Top = 115.27897277795276
$sh2.Width = 138.66016388031497
$sh2.Height = 33.92811023622047

$tr2 = $sh2.TextFrame.TextRange
$tr2.Text = "execute(“delete-patient n/John Doe”)"

# Recolor "execute" and the quoted command to accent1, leaving the
# punctuation/parentheses in the original blue (0070C0) -- this splits the
# single run into four runs matching the new wording.
$rExecute = $tr2.Characters(1, 7)
$rExecute.Font.Color.ObjectThemeColor = 5

$rCommand = $tr2.Characters(10, 25)
$rCommand.Font.Color.ObjectThemeColor = 5
